$wb = $excel.ActiveWorkbook

# Sheets "展览" (index 1) and "全部类型" (index 4) both hold the same table
# and both need: updated F2/F3/F5 counts, plus a new row inserted at row 6
# (pushing the old row 6 down to row 7, and bumping its "want to go" count).
$targetSheets = @(1, 4)

foreach ($idx in $targetSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # Bump the "想去人数" (want-to-go count) values for the existing rows.
    $ws.Range("F2").Value = 1284
    $ws.Range("F3").Value = 1619
    $ws.Range("F5").Value = 6206

    # Insert a new row before the current row 6 ("合肥·梦时空SPO1动漫展"),
    # which shifts that entire row down to row 7.
    $ws.Rows.Item(6).Insert()

    # The inserted row picked up a slightly different auto style; copy the
    # real column-A number style down from row 5 so A6 matches the rest.
    $ws.Range("A5").Copy()
    $ws.Range("A6").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    # Fill in the new event row. Force B6 to text first so the bare
    # "2024.04.04" string isn't auto-parsed into a date serial (matches how
    # the sibling date cells in this column are stored as text), then drop
    # back to the Normal style so no stray text-format style sticks around.
    $ws.Range("A6").Value = 5
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2024.04.04"
    $ws.Range("B6").Style = "Normal"
    $ws.Range("C6").Value = "合肥·环形宇宙动漫游戏嘉年华"
    $ws.Range("D6").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
    $ws.Range("E6").Value = "2024.04.04 09:30-04.05 17:00"
    $ws.Range("F6").Value = 12
    $ws.Range("G6").Value = "不可售"
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81916"
    $ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202402/1lGzmBT61708336972816.jpeg"

    # Renumber / bump the row that got pushed down to row 7.
    $ws.Range("A7").Value = 6
    $ws.Range("F7").Value = 106
}
